$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "25250"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1930"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2318"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "4355"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "3942"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "1780"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "9326"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "7192"
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "819"
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "5864"
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "852"
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "5224"
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2769"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "16299"
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "12522"
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "4732"
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "379"
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "6133"
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "6764"
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "6065"
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "3235"
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "4731"
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "858"
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "3867"
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "688"
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = "4520"
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "5307"
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = "2482"
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "4976"
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "4720"
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "3258"
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "6307"
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "17886"
$ws.Range("A37").NumberFormat = "@"
$ws.Range("A37").Value = "3776"
$ws.Range("A38").NumberFormat = "@"
$ws.Range("A38").Value = "5067"
$ws.Range("A39").NumberFormat = "@"
$ws.Range("A39").Value = "716"
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = "608"
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = "3315"
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "4340"
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "8347"
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = "5256"
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "31700"
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "2911"
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "11338"
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "1418"
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "50594"
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "48609"
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "47091"
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "5898"
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "19206"
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = "15610"
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "7470"
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "14787"
$ws.Range("A59").NumberFormat = "@"
$ws.Range("A59").Value = "7787"
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "32145"
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "10456"
$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "2588"
$ws.Range("A64").NumberFormat = "@"
$ws.Range("A64").Value = "42547"
$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = "17628"
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "20742"
$ws.Range("A69").NumberFormat = "@"
$ws.Range("A69").Value = "29988"
$ws.Range("A70").NumberFormat = "@"
$ws.Range("A70").Value = "18195"
$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "21276"
$ws.Range("A73").NumberFormat = "@"
$ws.Range("A73").Value = "8788"
$ws.Range("A74").NumberFormat = "@"
$ws.Range("A74").Value = "41217"
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "33492"
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "10093"
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = "10171"
$ws.Range("A80").NumberFormat = "@"
$ws.Range("A80").Value = "12479"
$ws.Range("A81").NumberFormat = "@"
$ws.Range("A81").Value = "12445"
$ws.Range("A82").NumberFormat = "@"
$ws.Range("A82").Value = "15860"
$ws.Range("A83").NumberFormat = "@"
$ws.Range("A83").Value = "49399"
$ws.Range("A84").NumberFormat = "@"
$ws.Range("A84").Value = "27135"
$ws.Range("A85").NumberFormat = "@"
$ws.Range("A85").Value = "31873"
$ws.Range("A87").NumberFormat = "@"
$ws.Range("A87").Value = "47616"
$ws.Range("A88").NumberFormat = "@"
$ws.Range("A88").Value = "37564"
$ws.Range("A89").NumberFormat = "@"
$ws.Range("A89").Value = "49748"
$ws.Range("A90").NumberFormat = "@"
$ws.Range("A90").Value = "49077"
$ws.Range("A92").NumberFormat = "@"
$ws.Range("A92").Value = "62830"
$ws.Range("A93").NumberFormat = "@"
$ws.Range("A93").Value = "43524"
$ws.Range("A94").NumberFormat = "@"
$ws.Range("A94").Value = "33053"
$ws.Range("A95").NumberFormat = "@"
$ws.Range("A95").Value = "62012"
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = "31842"
$ws.Range("A97").NumberFormat = "@"
$ws.Range("A97").Value = "47051"
$ws.Range("A99").NumberFormat = "@"
$ws.Range("A99").Value = "26599"
$ws.Range("A100").NumberFormat = "@"
$ws.Range("A100").Value = "52296"
$ws.Range("A101").NumberFormat = "@"
$ws.Range("A101").Value = "16085"
$ws.Range("A102").NumberFormat = "@"
$ws.Range("A102").Value = "24876"
$ws.Range("A103").NumberFormat = "@"
$ws.Range("A103").Value = "31152"
$ws.Range("A104").NumberFormat = "@"
$ws.Range("A104").Value = "17207"
$ws.Range("A105").NumberFormat = "@"
$ws.Range("A105").Value = "21673"
$ws.Range("A109").NumberFormat = "@"
$ws.Range("A109").Value = "26262"
$ws.Range("A110").NumberFormat = "@"
$ws.Range("A110").Value = "33608"
$ws.Range("A111").NumberFormat = "@"
$ws.Range("A111").Value = "54086"
$ws.Range("A112").NumberFormat = "@"
$ws.Range("A112").Value = "33663"
$ws.Range("A113").NumberFormat = "@"
$ws.Range("A113").Value = "38508"
$ws.Range("A114").NumberFormat = "@"
$ws.Range("A114").Value = "34751"
$ws.Range("A116").NumberFormat = "@"
$ws.Range("A116").Value = "75168"
$ws.Range("A117").NumberFormat = "@"
$ws.Range("A117").Value = "40180"
$ws.Range("A118").NumberFormat = "@"
$ws.Range("A118").Value = "29962"
$ws.Range("A120").NumberFormat = "@"
$ws.Range("A120").Value = "23229"
$ws.Range("A121").NumberFormat = "@"
$ws.Range("A121").Value = "63391"
$ws.Range("A122").NumberFormat = "@"
$ws.Range("A122").Value = "30307"
$ws.Range("A126").NumberFormat = "@"
$ws.Range("A126").Value = "31772"
$ws.Range("A131").NumberFormat = "@"
$ws.Range("A131").Value = "48417"
$ws.Range("A132").NumberFormat = "@"
$ws.Range("A132").Value = "61576"
$ws.Range("A133").NumberFormat = "@"
$ws.Range("A133").Value = "53526"
$ws.Range("A134").NumberFormat = "@"
$ws.Range("A134").Value = "44796"
$ws.Range("A135").NumberFormat = "@"
$ws.Range("A135").Value = "46991"
$ws.Range("A137").NumberFormat = "@"
$ws.Range("A137").Value = "50646"
$ws.Range("A138").NumberFormat = "@"
$ws.Range("A138").Value = "58054"
$ws.Range("A139").NumberFormat = "@"
$ws.Range("A139").Value = "63115"
$ws.Range("A140").NumberFormat = "@"
$ws.Range("A140").Value = "70263"
$ws.Range("A142").NumberFormat = "@"
$ws.Range("A142").Value = "45940"
$ws.Range("A143").NumberFormat = "@"
$ws.Range("A143").Value = "75264"
$ws.Range("A146").NumberFormat = "@"
$ws.Range("A146").Value = "31946"
$ws.Range("A148").NumberFormat = "@"
$ws.Range("A148").Value = "65570"
$ws.Range("A149").NumberFormat = "@"
$ws.Range("A149").Value = "6192"
